# "add property for block"
# Extend the XML-mapped table (表1) on Sheet1 with 12 new trailing columns:
#   Grass1..Grass5, Crack1..Crack5, Treasure1, Treasure2
# These describe new block properties (grass/crack/treasure decoration + tag slots).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newHeaders = @(
    "Grass1", "Grass2", "Grass3", "Grass4", "Grass5",
    "Crack1", "Crack2", "Crack3", "Crack4", "Crack5",
    "Treasure1", "Treasure2"
)

# 1) Grow the table (ListObject) so it spans the 12 extra columns (J:V instead of stopping at J).
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:V53"))

# 2) Re-assert the header text for every new column (must be written after Resize so the table
#    metadata actually keeps the real names instead of the generic "ColumnN" placeholders).
for ($i = 0; $i -lt $newHeaders.Count; $i++) {
    $ws.Cells.Item(1, 11 + $i).Value = $newHeaders[$i]
}

# 3) Widen the new columns to a readable size (matches the widths the authors left behind).
$ws.Columns.Item(11).ColumnWidth = 13.285714285714286   # K  Grass1
$ws.Columns.Item(12).ColumnWidth = 14.571428571428571   # L  Grass2
$ws.Columns.Item(13).ColumnWidth = 14.571428571428571   # M  Grass3
$ws.Columns.Item(14).ColumnWidth = 15.0                 # N  Grass4
$ws.Columns.Item(15).ColumnWidth = 13.142857142857142   # O  Grass5
$ws.Columns.Item(16).ColumnWidth = 16.428571428571427   # P  Crack1
$ws.Columns.Item(17).ColumnWidth = 9.857142857142858    # Q  Crack2
$ws.Columns.Item(18).ColumnWidth = 9.857142857142858    # R  Crack3
$ws.Columns.Item(19).ColumnWidth = 9.857142857142858    # S  Crack4
$ws.Columns.Item(21).ColumnWidth = 15.714285714285714   # U  Treasure1
$ws.Columns.Item(22).ColumnWidth = 14.714285714285714   # V  Treasure2

# 4) Move the selection/cursor to where the author ended up after adding the new fields.
$null = $ws.Range("U7").Select()
